$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 153.16667
$ws.Range("I9").Value = 173.66667
$ws.Range("J9").Value = 132.66667
$ws.Range("K9").Value = 173.66667
$ws.Range("L9").Value = 132.66667
$ws.Range("M9").Value = -4.666670000000011
$ws.Range("N9").Value = -470.66667
$ws.Range("H17").Value = 2291851.8
$ws.Range("J17").Value = 2902498.8
$ws.Range("L17").Value = 8707496.399999999
$ws.Range("N17").Value = -8707832.399999999
$ws.Range("H18").Value = 291.69232
$ws.Range("I18").Value = 291.69232
$ws.Range("K18").Value = 291.69232
$ws.Range("M18").Value = -7.692319999999995
$ws.Range("H80").Value = 772345.9399999999
$ws.Range("I80").Value = 2070.7144
$ws.Range("K80").Value = 6212.1432
$ws.Range("M80").Value = -5214.1432
$ws.Range("H83").Value = 772345.9399999999
$ws.Range("I83").Value = 2070.7144
$ws.Range("K83").Value = 18636.4296
$ws.Range("M83").Value = -13644.4296
$ws.Range("H99").Value = 622.2
$ws.Range("I99").Value = 587.6667
$ws.Range("K99").Value = 1763.0001
$ws.Range("M99").Value = -265.0001
$ws.Range("H118").Value = 857.1111
$ws.Range("I118").Value = 895
$ws.Range("K118").Value = 2685
$ws.Range("M118").Value = -1028
$ws.Range("H132").Value = 2000.75
$ws.Range("I132").Value = 2148.8572
$ws.Range("K132").Value = 6446.571599999999
$ws.Range("M132").Value = -3916.571599999999
$ws.Range("H137").Value = 20517.52
$ws.Range("J137").Value = 24791.176
$ws.Range("L137").Value = 74373.528
$ws.Range("N137").Value = -79473.528
$ws.Range("H138").Value = 1890.4546
$ws.Range("I138").Value = 1425.1818
$ws.Range("J138").Value = 2355.7273
$ws.Range("K138").Value = 4275.5454
$ws.Range("L138").Value = 7067.1819
$ws.Range("M138").Value = 864.4546
$ws.Range("N138").Value = -17347.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4411.108
$ws.Range("I32").Value = 2865.0986
$ws.Range("J32").Value = 41000
$ws.Range("K32").Value = 2865.0986
$ws.Range("L32").Value = 41000
$ws.Range("M32").Value = -2578.0986
$ws.Range("N32").Value = -41574
$ws.Range("H45").Value = 9865.75
$ws.Range("I45").Value = 10393.467
$ws.Range("K45").Value = 10393.467
$ws.Range("M45").Value = -10016.467
$ws.Range("H61").Value = 4118.0444
$ws.Range("I61").Value = 2953.025
$ws.Range("K61").Value = 2953.025
$ws.Range("M61").Value = -2741.025
$ws.Range("H122").Value = 1799.4546
$ws.Range("J122").Value = 3142.8572
$ws.Range("L122").Value = 9428.571599999999
$ws.Range("N122").Value = -14328.5716
$ws.Range("H136").Value = 4118.0444
$ws.Range("I136").Value = 2953.025
$ws.Range("K136").Value = 8859.075000000001
$ws.Range("M136").Value = -6309.075000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3430.75
$ws.Range("I105").Value = 3662.5715
$ws.Range("K105").Value = 3662.5715
$ws.Range("M105").Value = -1915.5715
$ws.Range("H134").Value = 12355.448
$ws.Range("I134").Value = 6231.421
$ws.Range("J134").Value = 23991.1
$ws.Range("K134").Value = 18694.263
$ws.Range("L134").Value = 71973.29999999999
$ws.Range("M134").Value = -16159.263
$ws.Range("N134").Value = -77043.29999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 337003.66
$ws.Range("I31").Value = 58696.832
$ws.Range("J31").Value = 670971.9
$ws.Range("K31").Value = 58696.832
$ws.Range("L31").Value = 670971.9
$ws.Range("M31").Value = -58401.832
$ws.Range("N31").Value = -671561.9
$ws.Range("H34").Value = 337003.66
$ws.Range("I34").Value = 58696.832
$ws.Range("J34").Value = 670971.9
$ws.Range("K34").Value = 58696.832
$ws.Range("L34").Value = 670971.9
$ws.Range("M34").Value = -58494.832
$ws.Range("N34").Value = -671375.9
$ws.Range("H41").Value = 37551.332
$ws.Range("I41").Value = 18999.5
$ws.Range("K41").Value = 18999.5
$ws.Range("M41").Value = -18571.5
$ws.Range("H58").Value = 3763.7144
$ws.Range("I58").Value = 2033.9474
$ws.Range("K58").Value = 2033.9474
$ws.Range("M58").Value = -1830.9474
$ws.Range("H134").Value = 3439.5312
$ws.Range("I134").Value = 2719.724
$ws.Range("J134").Value = 10397.667
$ws.Range("K134").Value = 8159.172
$ws.Range("L134").Value = 31193.001
$ws.Range("M134").Value = -5624.172
$ws.Range("N134").Value = -36263.001
$ws.Range("H136").Value = 3763.7144
$ws.Range("I136").Value = 2033.9474
$ws.Range("K136").Value = 6101.8422
$ws.Range("M136").Value = -3551.8422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1985.0588
$ws.Range("I5").Value = 716
$ws.Range("J5").Value = 2257
$ws.Range("K5").Value = 2148
$ws.Range("L5").Value = 6771
$ws.Range("M5").Value = -2036
$ws.Range("N5").Value = -6995
$ws.Range("H57").Value = 1500
$ws.Range("I57").Value = 1500
$ws.Range("K57").Value = 4500
$ws.Range("M57").Value = -3941
$ws.Range("H113").Value = 1088.2941
$ws.Range("I113").Value = 1239.5714
$ws.Range("J113").Value = 982.4
$ws.Range("K113").Value = 3718.7142
$ws.Range("L113").Value = 2947.2
$ws.Range("M113").Value = -1548.7142
$ws.Range("N113").Value = -7287.2
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H135").Value = 1985.0588
$ws.Range("I135").Value = 716
$ws.Range("J135").Value = 2257
$ws.Range("K135").Value = 6444
$ws.Range("L135").Value = 20313
$ws.Range("M135").Value = -3909
$ws.Range("N135").Value = -25383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 14287945
$ws.Range("I102").Value = 2479.111
$ws.Range("K102").Value = 2479.111
$ws.Range("M102").Value = -857.1109999999999
$ws.Range("H107").Value = 313.77777
$ws.Range("I107").Value = 335.55
$ws.Range("J107").Value = 251.57143
$ws.Range("K107").Value = 335.55
$ws.Range("L107").Value = 251.57143
$ws.Range("M107").Value = 1584.45
$ws.Range("N107").Value = -4091.57143
$ws.Range("H122").Value = 5931.492
$ws.Range("I122").Value = 4268.3
$ws.Range("J122").Value = 8592.6
$ws.Range("K122").Value = 12804.9
$ws.Range("L122").Value = 25777.8
$ws.Range("M122").Value = -10354.9
$ws.Range("N122").Value = -30677.8
$ws.Range("H126").Value = 4353.9473
$ws.Range("I126").Value = 4464.25
$ws.Range("J126").Value = 3765.6667
$ws.Range("K126").Value = 13392.75
$ws.Range("L126").Value = 11297.0001
$ws.Range("M126").Value = -10922.75
$ws.Range("N126").Value = -16237.0001
$ws.Range("H132").Value = 26285.732
$ws.Range("I132").Value = 18175.924
$ws.Range("K132").Value = 54527.772
$ws.Range("M132").Value = -51997.772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3622.7144
$ws.Range("J100").Value = 3998.7693
$ws.Range("L100").Value = 3998.7693
$ws.Range("N100").Value = -5080.7693
$ws.Range("H122").Value = 5331.3335
$ws.Range("J122").Value = 4997.25
$ws.Range("L122").Value = 14991.75
$ws.Range("N122").Value = -19891.75
$ws.Range("H132").Value = 2006.2787
$ws.Range("I132").Value = 1250.3778
$ws.Range("J132").Value = 4132.25
$ws.Range("K132").Value = 3751.1334
$ws.Range("L132").Value = 12396.75
$ws.Range("M132").Value = -1221.1334
$ws.Range("N132").Value = -17456.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28337.174
$ws.Range("I132").Value = 27785.246
$ws.Range("K132").Value = 83355.738
$ws.Range("M132").Value = -80825.738
